$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.222.29"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.859.25"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +1.33%  "
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "311.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3723"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07311"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9361"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07871"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "1.869.62"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.531"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.022"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008728"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "27.239.25"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.103"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.973"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08884"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.344"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.583"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7406"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.684"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02025"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5322"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.099"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.332"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4784"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.019"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.634"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "66.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06076"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8989"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.17%  "
